$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header / summary cell updates
# ---------------------------------------------------------------------------
$ws.Range("D5").Value = "Report Generated On: 08/18/2025 09:47 PM"
$ws.Range("C8").Value = 0
$ws.Range("C9").Value = 15

# ---------------------------------------------------------------------------
# 2. Thursday block (rows 16-18): pricing zeroed out
# ---------------------------------------------------------------------------
$ws.Range("H16").Value = 0
$ws.Range("H17").Value = 0
$ws.Range("H18").Value = 0

# ---------------------------------------------------------------------------
# 3. Friday block: insert 4 new "TIE-4-ALH-F" line items, one after each
#    Point group (01, 03, 05, 07), and zero-out every pricing (H) cell.
#    Work bottom-to-top (using the ORIGINAL row numbers) so earlier
#    insertion points never move while we still need them.
# ---------------------------------------------------------------------------
$ws.Rows.Item(32).Insert()
$ws.Rows.Item(30).Insert()
$ws.Rows.Item(28).Insert()
$ws.Rows.Item(26).Insert()

# After the four inserts above, the sheet rows 23-36 line up as:
#   23 Point01 CON-10-AAA-1-B      (original row 23, untouched)
#   24 Point01 CON-10-AAA-3-P      (original row 24, untouched)
#   25 Point01 CON-40-AAA-1-B      (original row 25, untouched)
#   26 Point01 TIE-4-ALH-F         NEW (gray stripe)
#   27 Point03 CON-10-AAA-1-B-REEL (was original row 26)
#   28 Point03 CON-40-AAA-1-B      (was original row 27)
#   29 Point03 TIE-4-ALH-F         NEW (white stripe)
#   30 Point05 CON-10-AAA-1-B-REEL (was original row 28)
#   31 Point05 CON-40-AAA-1-B      (was original row 29)
#   32 Point05 TIE-4-ALH-F         NEW (gray stripe)
#   33 Point07 CON-10-AAA-1-B-REEL (was original row 30)
#   34 Point07 CON-40-AAA-1-B      (was original row 31)
#   35 Point07 TIE-4-ALH-F         NEW (white stripe)
#   36 TOTAL                       (was original row 32)

# Make sure the four brand-new rows carry the correct zebra-stripe style
# (gray fill = style used by row 24, white = style used by row 23) by
# pasting formats from those untouched reference rows.
$ws.Range("A24:I24").Copy()
$ws.Range("A26:I26").PasteSpecial(-4122)
$ws.Range("A32:I32").PasteSpecial(-4122)

$ws.Range("A23:I23").Copy()
$ws.Range("A29:I29").PasteSpecial(-4122)
$ws.Range("A35:I35").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4. Fill in the content for the four new TIE-4-ALH-F rows.
# ---------------------------------------------------------------------------
$ws.Range("A26").Value = "Point 01"
$ws.Range("B26").Value = "TIE-4-ALH-F"
$ws.Range("C26").Value = "Inst"
$ws.Range("D26").Value = "TIE,4 AWG,AL Hand Tie,F Neck"
$ws.Range("E26").Value = "EA"
$ws.Range("F26").Value = 12
$ws.Range("H26").Value = 0

$ws.Range("A29").Value = "Point 03"
$ws.Range("B29").Value = "TIE-4-ALH-F"
$ws.Range("C29").Value = "Inst"
$ws.Range("D29").Value = "TIE,4 AWG,AL Hand Tie,F Neck"
$ws.Range("E29").Value = "EA"
$ws.Range("F29").Value = 18
$ws.Range("H29").Value = 0

$ws.Range("A32").Value = "Point 05"
$ws.Range("B32").Value = "TIE-4-ALH-F"
$ws.Range("C32").Value = "Inst"
$ws.Range("D32").Value = "TIE,4 AWG,AL Hand Tie,F Neck"
$ws.Range("E32").Value = "EA"
$ws.Range("F32").Value = 18
$ws.Range("H32").Value = 0

$ws.Range("A35").Value = "Point 07"
$ws.Range("B35").Value = "TIE-4-ALH-F"
$ws.Range("C35").Value = "Inst"
$ws.Range("D35").Value = "TIE,4 AWG,AL Hand Tie,F Neck"
$ws.Range("E35").Value = "EA"
$ws.Range("F35").Value = 24
$ws.Range("H35").Value = 0

# ---------------------------------------------------------------------------
# 5. Zero-out the pricing column for every pre-existing Friday-block row
#    (the quantities/descriptions themselves are unchanged).
# ---------------------------------------------------------------------------
$ws.Range("H23").Value = 0
$ws.Range("H24").Value = 0
$ws.Range("H25").Value = 0
$ws.Range("H27").Value = 0
$ws.Range("H28").Value = 0
$ws.Range("H30").Value = 0
$ws.Range("H31").Value = 0
$ws.Range("H33").Value = 0
$ws.Range("H34").Value = 0

# ---------------------------------------------------------------------------
# 6. TOTAL row (now row 36) pricing.
# ---------------------------------------------------------------------------
$ws.Range("H36").Value = 0

Write-Host "Edit complete"
